# Add ten new device rows (Finger Print Scanner 30/31, IRIS Scanner 30/31,
# Web Camera 30/31, Document Scanner 30/31, Printer 30/31) to the bottom of
# the device master data table, continuing the existing pattern of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 147

$rows = @(
    @(3000166, "Finger Print Scanner 30", "D6-15-AC-80-6B-86", "BS563Q2230814", 165),
    @(3000167, "IRIS Scanner 30",         "6D-58-E2-DF-74-34", "BS563Q2230815", 327),
    @(3000168, "Web Camera 30",           "E2-A8-56-86-15-30", "BS563Q2230816", 736),
    @(3000169, "Document Scanner 30",     "72-E8-B9-FD-63-65", "BS563Q2230817", 801),
    @(3000170, "Printer 30",              "D3-F3-A4-50-AD-12", "BS563Q2230818", 920),
    @(3000171, "Finger Print Scanner 31", "06-16-D0-0B-A6-E4", "BS563Q2230819", 165),
    @(3000172, "IRIS Scanner 31",         "21-78-45-AC-E9-20", "BS563Q2230820", 327),
    @(3000173, "Web Camera 31",           "3C-E8-87-99-DB-FA", "BS563Q2230821", 736),
    @(3000174, "Document Scanner 31",     "BF-55-53-98-40-08", "BS563Q2230822", 801),
    @(3000175, "Printer 31",              "5A-43-36-46-22-EB", "BS563Q2230823", 920)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]          # id
    $ws.Cells.Item($r, 2).Value = $data[1]          # name
    $ws.Cells.Item($r, 3).Value = $data[2]          # mac_address
    $ws.Cells.Item($r, 4).Value = $data[3]          # serial_num
    $ws.Cells.Item($r, 6).Value = $data[4]          # dspec_id
    $ws.Cells.Item($r, 7).Value = "eng"             # lang_code
    $ws.Cells.Item($r, 8).Value = $true             # is_active
    $ws.Cells.Item($r, 8).HorizontalAlignment = -4131   # xlHAlignLeft (matches existing rows' style)
    $ws.Cells.Item($r, 9).Value = "superadmin"      # cr_by
    $ws.Cells.Item($r, 10).Value = "now()"          # cr_dtimes
    $ws.Cells.Item($r, 11).Value = "now()"          # eff_dtimes
}

$lastRow = $startRow + $rows.Count - 1

# Match the author's final view state: selection on D145, scrolled so row 139
# is the top visible row.
$ws.Range("D145").Select()
$win = $excel.Windows.Item(1)
$win.ScrollRow = 139
$win.ScrollColumn = 1
